# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.207.44"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "'1.681.16"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'216.34"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "'0.5249"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").Value = "'21.44"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("D11").Value = "'0.07618"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "'1.705.12"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "'4.517"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "'0.5742"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "'0.000008244"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "'66.08"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").Value = "'26.234.90"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'4.870"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'10.76"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "'189.30"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D23").Value = "'1.007"
$ws.Range("D24").Value = "'148.86"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").Value = "'0.1259"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'7.730"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "'15.78"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "'0.06376"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "'1.378"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "'1.315"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "'3.568"
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").Value = "'1.680"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").Value = "'1.022"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "'0.6111"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").Value = "'2.423"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "'2.748"
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("D38").Value = "'6.169"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").Value = "'0.01613"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'1.097.04"
$ws.Range("D41").Value = "'0.8847"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'100.50"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'1.832.25"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").Value = "'57.37"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "'0.05270"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'0.4280"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "'6.000"
$ws.Range("E51").Value = "  -1.16%  "
